# PM Review Slides - minor edits
# The "Iteration 2" date-legend strip on slide 12 (the Gantt-like date boxes)
# shifts forward by one day: each box now shows the date previously shown by
# the next box, "3/10" is dropped (replaced by "4/10"), and the final box
# becomes "14/10 -  15/10" instead of "13/10 -  15/10".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$updates = @{
    76 = "4/10"
    77 = "5/10"
    78 = "6/10"
    79 = "7/10"
    80 = "8/10"
    81 = "9/10"
    82 = "10/10"
    83 = "11/10"
    84 = "12/10"
    85 = "13/10"
}

foreach ($id in $updates.Keys) {
    $shp = Get-ShapeById $s $id
    $shp.TextFrame.TextRange.Text = $updates[$id]
}

# Last box keeps the "-  15/10" tail but the first number moves from 13 to 14,
# split across two runs: "14/10 " and "-  15/10".
$lastShp = Get-ShapeById $s 86
$tr = $lastShp.TextFrame.TextRange
$firstPart = $tr.Characters(1, 6)
$firstPart.Text = "14/10 "
